# durham_daily_vaccinations.xlsx update
# - Row 110 (2021-04-09) loses its "latest row" highlight (Good/green) and
#   reverts to the regular alternating "Neutral"/yellow look, plus a fixed
#   15pt row height. Its "daily rate to achieve June 20 target" label moves
#   down to the new last row.
# - A brand new data row 111 (2021-04-10, 3984 doses) is appended, taking
#   over the green "Good" highlight + the June-20-target label.
# - The old two scratch/summary rows (113: SUM(B2:B102), 114: the target
#   delta) are removed, and a couple of blank spacer rows are left below
#   the data instead (115 and 117, with 116 skipped entirely, mirroring a
#   manual row-113/114 content clear + a couple of stray blank rows typed
#   further down).
# - Selection moves to H120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$goodFont = 24832       # RGB(0,97,0)   -> "Good" cell style font color
$goodFill = 13561798    # RGB(198,239,206) -> "Good" cell style fill color
$neutralFont = 22428    # RGB(156,87,0) -> "Neutral" cell style font color
$neutralFill = 10284031 # RGB(255,235,156) -> "Neutral" cell style fill color

# ---------------------------------------------------------------------
# Row 110: Good (green) -> Neutral (yellow), fixed row height, drop I110
# ---------------------------------------------------------------------
$ws.Range("A110").Font.Color = $neutralFont
$ws.Range("A110").Interior.Color = $neutralFill
$ws.Range("A110").NumberFormat = "d-mmm"

$ws.Range("B110").Font.Color = $neutralFont
$ws.Range("B110").Interior.Color = $neutralFill

$ws.Range("C110").Font.Color = $neutralFont
$ws.Range("C110").Interior.Color = $neutralFill
$ws.Range("C110").NumberFormat = "0"

$ws.Range("D110").Font.Color = $neutralFont
$ws.Range("D110").Interior.Color = $neutralFill

$ws.Range("E110").Font.Color = $neutralFont
$ws.Range("E110").Interior.Color = $neutralFill

$ws.Range("F110").Font.Color = $neutralFont
$ws.Range("F110").Interior.Color = $neutralFill
$ws.Range("F110").NumberFormat = "0"

$ws.Range("G110").Font.Color = $neutralFont
$ws.Range("G110").Interior.Color = $neutralFill

$ws.Range("H110").Font.Color = $neutralFont
$ws.Range("H110").Interior.Color = $neutralFill
$ws.Range("H110").NumberFormat = "0"

$ws.Range("I110").Clear()

$ws.Rows.Item(110).RowHeight = 15

# ---------------------------------------------------------------------
# Row 111 (new): the new day's data, styled Good (green), carries the
# "daily rate to achieve June 20 target" label that used to sit on 110.
# ---------------------------------------------------------------------
$ws.Range("A111").Value = 44296
$ws.Range("A111").Font.Color = $goodFont
$ws.Range("A111").Interior.Color = $goodFill
$ws.Range("A111").NumberFormat = "d-mmm"

$ws.Range("B111").Value = 3984
$ws.Range("B111").Font.Color = $goodFont
$ws.Range("B111").Interior.Color = $goodFill

$ws.Range("C111").Formula = "=(AVERAGE(B105:B111))"
$ws.Range("C111").Font.Color = $goodFont
$ws.Range("C111").Interior.Color = $goodFill
$ws.Range("C111").NumberFormat = "0"

$ws.Range("D111").Formula = "=(D110-B111)"
$ws.Range("D111").Font.Color = $goodFont
$ws.Range("D111").Interior.Color = $goodFill

$ws.Range("E111").Formula = "=E110+B111"
$ws.Range("E111").Font.Color = $goodFont
$ws.Range("E111").Interior.Color = $goodFill

$ws.Range("F111").Formula = "=D111/C111"
$ws.Range("F111").Font.Color = $goodFont
$ws.Range("F111").Interior.Color = $goodFill
$ws.Range("F111").NumberFormat = "0"

$ws.Range("G111").Formula = "=A111+F111"
$ws.Range("G111").Font.Color = $goodFont
$ws.Range("G111").Interior.Color = $goodFill

$ws.Range("H111").Formula = "=D111/84"
$ws.Range("H111").Font.Color = $goodFont
$ws.Range("H111").Interior.Color = $goodFill
$ws.Range("H111").NumberFormat = "0"

$ws.Range("I111").Value = "daily rate to achieve June 20 target"
$ws.Range("I111").Font.Color = $goodFont
$ws.Range("I111").Interior.Color = $goodFill

# ---------------------------------------------------------------------
# Remove the old scratch rows (113: SUM(B2:B102), 114: the target delta)
# and leave a couple of blank spacer rows below the data (115, 117 - 116
# stays untouched/absent).
# ---------------------------------------------------------------------
$ws.Range("A113:B113").Clear()
$ws.Range("B114").Clear()

$ws.Range("A115").NumberFormat = "d-mmm"
$ws.Range("A117").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------
# Selection moves to H120
# ---------------------------------------------------------------------
$ws.Range("H120").Select() | Out-Null
